$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "yes" answer and the feature/java test references for rows 19-27
# (HTB royal flush beats AIP <hand>)
for ($row = 19; $row -le 27; $row++) {
    $ws.Range("B$row").Value = "yes"
    $ws.Range("E$row").Value = "HTB_royal_flush_beats_AIP.feature"
    $ws.Range("F$row").Value = "HTBRoyalFlushBeatsAIP.java"
}

# Update the selection / scroll position to reflect where the user ended up
$ws.Range("F28").Select() | Out-Null
